# Update the skoda_octavia mileage data: refresh existing year rows (2001-2013)
# with newly scraped values and append rows for the newly scraped years
# 2014-2018, then extend the two chart series ranges to cover the larger
# data block.

$wb = $excel.ActiveWorkbook
$sheetName = "skoda_octavia"
$ws = $null
try {
    $ws = $wb.Worksheets.Item($sheetName)
} catch {
    $ws = $null
}
if ($ws -eq $null) {
    $ws = $wb.Worksheets.Add()
    $ws.Name = $sheetName
}

# row, year label (note the trailing space - matches the source data), avg. price (B), running/overall avg (C)
$data = @(
    @(2,  "2001 ", 175000,    178028.3333333333),
    @(3,  "2002 ", 198000,    198000),
    @(4,  "2003 ", 145189,    145189),
    @(5,  "2005 ", 143794,    176414.3333333333),
    @(6,  "2006 ", 197000,    192637.6666666667),
    @(7,  "2007 ", 198000,    208333.3333333333),
    @(8,  "2008 ", 181322.5,  183912),
    @(9,  "2009 ", 136125.5,  140462.25),
    @(10, "2010 ", 102317,    115111.3333333333),
    @(11, "2011 ", 185500,    174475),
    @(12, "2012 ", 118167.5,  118387.375),
    @(13, "2013 ", 130000,    132304.9090909091),
    @(14, "2014 ", 143057.5,  140696.05),
    @(15, "2015 ", 82000,     105453),
    @(16, "2016 ", 27200,     27515.28571428571),
    @(17, "2017 ", 1,         2584.916666666667),
    @(18, "2018 ", 1,         1)
)

foreach ($row in $data) {
    $r = $row[0]
    $yearLabel = $row[1]
    $avgPrice = $row[2]
    $runningAvg = $row[3]

    $aCell = $ws.Cells.Item($r, 1)
    # Prefix with an apostrophe so the numeric-looking text ("2001 ") is
    # stored as text instead of being coerced into a number, then strip the
    # formatting stamp COM leaves behind so the cell stays plain text.
    $aCell.Value = "'" + $yearLabel
    $aCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $avgPrice
    $ws.Cells.Item($r, 3).Value = $runningAvg
}

# Extend the two bar-chart series (avg price / running avg) to the new
# data range (still offset by one row like the original chart), if the
# sheet already has a chart on it.
$lastRow = $data[$data.Length - 1][0]
$chartLastRow = $lastRow - 1
if ($ws.ChartObjects().Count -gt 0) {
    $co = $ws.ChartObjects().Item(1)
    $chart = $co.Chart
    $series = $chart.SeriesCollection()
    $series.Item(1).Formula = "=SERIES(,skoda_octavia!`$A`$1:`$A`$$chartLastRow,skoda_octavia!`$B`$1:`$B`$$chartLastRow,1)"
    $series.Item(2).Formula = "=SERIES(,skoda_octavia!`$A`$1:`$A`$$chartLastRow,skoda_octavia!`$C`$1:`$C`$$chartLastRow,2)"
}
